$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.899.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.288.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.91"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.630"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.07"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +5.47%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.87%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.91%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.49"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +16.26%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.627.19"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.82"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.94"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.65%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.290.35"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.784.27"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0947"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.32"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.19"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.98"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +8.19%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.86"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.47"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.139"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.57"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.02%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.81%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0703"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.17"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.74"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.70"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0248"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.98%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.12"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +28.56%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.67"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.60%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.53"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000219"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -8.08%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0964"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.09"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.488.77"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.95"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.36"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.90%  "
